$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45
$ws.Range("B45").Value = "Scaler"
$ws.Range("A45").Value = "Tree"
$ws.Range("C45").Value = @"
Given two binary trees, check if they are equal or not.
Two binary trees are considered equal if they are structurally identical and the nodes have the same value.
"@
$ws.Range("C45").WrapText = $true
$ws.Range("D45").Value = @"
Solution39
"@
$ws.Rows.Item(45).RowHeight = 60

# Row 46
$ws.Range("B46").Value = "Scaler"
$ws.Range("A46").Value = "Array"
$ws.Range("C46").Value = @"
Given an array of integers A, a subarray of an array is said to be good if it fulfills any one of the criteria:
1. Length of the subarray is be even, and the sum of all the elements of the subarray must be less than B.
2. Length of the subarray is be odd, and the sum of all the elements of the subarray must be greater than B.
Your task is to find the count of good subarrays in A.
"@
$ws.Range("C46").Font.Name = "Arial"
$ws.Range("C46").Font.Size = 12
$ws.Range("C46").Font.Color = 5329233
$ws.Range("C46").WrapText = $true
$ws.Range("D46").Value = @"
Solution40
"@
$ws.Rows.Item(46).RowHeight = 105.75

# Row 47
$ws.Range("B47").Value = "LeetCode"
$ws.Range("A47").Value = "String"
$ws.Range("C47").Value = @"
iven a string s, find the length of the longest substring without repeating characters.
"@
$ws.Range("C47").Font.Name = "Arial"
$ws.Range("C47").Font.Size = 12
$ws.Range("C47").Font.Color = 5329233
$ws.Range("C47").WrapText = $true
$ws.Range("D47").Value = @"
Solution41
"@
$ws.Rows.Item(47).RowHeight = 30.75

# Row 48
$ws.Range("B48").Value = "Scaler"
$ws.Range("A48").Value = "Tree"
$ws.Range("C48").Value = @"
Given a binary tree, check whether it is a mirror of itself (i.e., symmetric around its center).
"@
$ws.Range("C48").Font.Name = "Arial"
$ws.Range("C48").Font.Size = 12
$ws.Range("C48").Font.Color = 5329233
$ws.Range("C48").WrapText = $true
$ws.Range("D48").Value = @"
Solution42
"@
$ws.Rows.Item(48).RowHeight = 30.75

# Row 49
$ws.Range("B49").Value = "Scaler"
$ws.Range("C49").Value = @"
Given a set of distinct integers A, return all possible subsets.
NOTE:
Elements in a subset must be in non-descending order.
The solution set must not contain duplicate subsets.
Also, the subsets should be sorted in ascending ( lexicographic ) order.
The list is not necessarily sorted.
"@
$ws.Range("C49").Font.Name = "Arial"
$ws.Range("C49").Font.Size = 12
$ws.Range("C49").Font.Color = 5329233
$ws.Range("C49").WrapText = $true
$ws.Range("D49").Value = @"
Solution43
"@
$ws.Range("A49").Value = "Subsequence"
$ws.Rows.Item(49).RowHeight = 120.75

# Row 50
$ws.Range("B50").Value = "Scaler"
$ws.Range("A50").Value = "Tree"
$ws.Range("C50").Value = @"
Given a binary tree A, invert the binary tree and return it.
Inverting refers to making the left child the right child and vice versa.
"@
$ws.Range("C50").WrapText = $true
$ws.Range("D50").Value = @"
Solution44
"@
$ws.Rows.Item(50).RowHeight = 45

# Row 51
$ws.Range("B51").Value = "Scaler"
$ws.Range("A51").Value = "Tree"
$ws.Range("D51").Value = @"
Solution45
"@
$ws.Range("C51").Value = @"
Given a Binary Tree A containing N nodes, you need to find the path from Root to a given node B.
"@
$ws.Range("C51").Font.Name = "Arial"
$ws.Range("C51").Font.Size = 12
$ws.Range("C51").Font.Color = 5329233
$ws.Range("C51").WrapText = $true
$ws.Rows.Item(51).RowHeight = 30.75

# Row 52
$ws.Range("B52").Value = "Scaler"
$ws.Range("C52").Value = @"
You are given an array of integers of N size.
You have to find that there is any subsequence exists or not whose sum is equal to B.
"@
$ws.Range("C52").WrapText = $true
$ws.Range("D52").Value = @"
Solution46
"@
$ws.Range("A52").Value = "Subsequence"
$ws.Rows.Item(52).RowHeight = 45

# Update data validation list for column A to include Subsequence
$v = $ws.Range("A1").Validation
$v.Formula1 = '"Subsequence,Tree,Stack,Queue,Recursion,Hashing,Array, String,Bit Manupulation,Loop,Maths,Modulus"'

# Update view state to match post-edit selection
$ws.Range("C51").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1
